$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - reuse H1's header style
# (bold font, thin border, centered) by copying its format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2 (plain numeric values, no special style)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
